$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 215.12
$ws.Range("I15").Value = 215.12
$ws.Range("K15").Value = 645.36
$ws.Range("M15").Value = -476.36
$ws.Range("H28").Value = 819.1
$ws.Range("I28").Value = 951.7646999999999
$ws.Range("K28").Value = 951.7646999999999
$ws.Range("M28").Value = -466.7646999999999
$ws.Range("H113").Value = 6073
$ws.Range("I113").Value = 6168.3335
$ws.Range("J113").Value = 6001.5
$ws.Range("K113").Value = 6168.3335
$ws.Range("L113").Value = 6001.5
$ws.Range("M113").Value = -2914.3335
$ws.Range("N113").Value = -12509.5
$ws.Range("H132").Value = 2085200.1
$ws.Range("I132").Value = 1811.881
$ws.Range("J132").Value = 16668918
$ws.Range("K132").Value = 5435.643
$ws.Range("L132").Value = 50006754
$ws.Range("M132").Value = -2905.643
$ws.Range("N132").Value = -50011814
$ws.Range("H138").Value = 2619.7273
$ws.Range("I138").Value = 883.91113
$ws.Range("J138").Value = 4066.2407
$ws.Range("K138").Value = 2651.73339
$ws.Range("L138").Value = 12198.7221
$ws.Range("M138").Value = 2488.26661
$ws.Range("N138").Value = -22478.7221

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 5601.9375
$ws.Range("I102").Value = 4158.5713
$ws.Range("J102").Value = 15705.5
$ws.Range("K102").Value = 4158.5713
$ws.Range("L102").Value = 15705.5
$ws.Range("M102").Value = -2536.5713
$ws.Range("N102").Value = -18949.5
$ws.Range("H132").Value = 1135.898
$ws.Range("I132").Value = 1045.4166
$ws.Range("J132").Value = 1386.4615
$ws.Range("K132").Value = 3136.2498
$ws.Range("L132").Value = 4159.3845
$ws.Range("M132").Value = -606.2498000000001
$ws.Range("N132").Value = -9219.3845

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 83335880
$ws.Range("J99").Value = 1533.3334
$ws.Range("L99").Value = 1533.3334
$ws.Range("N99").Value = -4529.3334
$ws.Range("H105").Value = 4107.3555
$ws.Range("I105").Value = 3768.3872
$ws.Range("J105").Value = 4857.9287
$ws.Range("K105").Value = 3768.3872
$ws.Range("L105").Value = 4857.9287
$ws.Range("M105").Value = -2021.3872
$ws.Range("N105").Value = -8351.9287

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 45702.523
$ws.Range("I31").Value = 2003.625
$ws.Range("J31").Value = 145585.72
$ws.Range("K31").Value = 2003.625
$ws.Range("L31").Value = 145585.72
$ws.Range("M31").Value = -1708.625
$ws.Range("N31").Value = -146175.72
$ws.Range("H34").Value = 45702.523
$ws.Range("I34").Value = 2003.625
$ws.Range("J34").Value = 145585.72
$ws.Range("K34").Value = 2003.625
$ws.Range("L34").Value = 145585.72
$ws.Range("M34").Value = -1801.625
$ws.Range("N34").Value = -145989.72
$ws.Range("H86").Value = 2146.3333
$ws.Range("I86").Value = 2174.8
$ws.Range("J86").Value = 2004
$ws.Range("K86").Value = 2174.8
$ws.Range("L86").Value = 2004
$ws.Range("M86").Value = -1051.8
$ws.Range("N86").Value = -4250
$ws.Range("H89").Value = 2146.3333
$ws.Range("I89").Value = 2174.8
$ws.Range("J89").Value = 2004
$ws.Range("K89").Value = 10874
$ws.Range("L89").Value = 10020
$ws.Range("M89").Value = -5258
$ws.Range("N89").Value = -21252
$ws.Range("H99").Value = 3508.6956
$ws.Range("I99").Value = 2308.3333
$ws.Range("K99").Value = 2308.3333
$ws.Range("M99").Value = -810.3332999999998
$ws.Range("H126").Value = 3508.6956
$ws.Range("I126").Value = 2308.3333
$ws.Range("K126").Value = 6924.999899999999
$ws.Range("M126").Value = -4454.999899999999
$ws.Range("H134").Value = 16667875
$ws.Range("I134").Value = 1066.8846
$ws.Range("J134").Value = 125002130
$ws.Range("K134").Value = 3200.6538
$ws.Range("L134").Value = 375006390
$ws.Range("M134").Value = -665.6538
$ws.Range("N134").Value = -375011460

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 16194847
$ws.Range("I131").Value = 83500300
$ws.Range("J131").Value = 41536.96
$ws.Range("K131").Value = 250500900
$ws.Range("L131").Value = 124610.88
$ws.Range("M131").Value = -250495860
$ws.Range("N131").Value = -134690.88
$ws.Range("H132").Value = 816.5263
$ws.Range("I132").Value = 784.35297
$ws.Range("J132").Value = 1090
$ws.Range("K132").Value = 7059.17673
$ws.Range("L132").Value = 9810
$ws.Range("M132").Value = -4529.17673
$ws.Range("N132").Value = -14870

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1300
$ws.Range("I122").Value = 1300
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3900
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1450
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 2516.6897
$ws.Range("I132").Value = 2082.7917
$ws.Range("J132").Value = 4599.4
$ws.Range("K132").Value = 6248.375100000001
$ws.Range("L132").Value = 13798.2
$ws.Range("M132").Value = -3718.375100000001
$ws.Range("N132").Value = -18858.2
$ws.Range("H134").Value = 20853.25
$ws.Range("J134").Value = 20853.25
$ws.Range("L134").Value = 62559.75
$ws.Range("N134").Value = -67629.75
$ws.Range("H136").Value = 12025.895
$ws.Range("J136").Value = 12025.895
$ws.Range("L136").Value = 36077.685
$ws.Range("N136").Value = -41177.685

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2650
$ws.Range("I40").Value = 2100
$ws.Range("J40").Value = 3200
$ws.Range("K40").Value = 2100
$ws.Range("L40").Value = 3200
$ws.Range("M40").Value = -1964
$ws.Range("N40").Value = -3472
$ws.Range("H93").Value = 3014.9
$ws.Range("I93").Value = 2935.4285
$ws.Range("J93").Value = 3057.6924
$ws.Range("K93").Value = 2935.4285
$ws.Range("L93").Value = 3057.6924
$ws.Range("M93").Value = -1687.4285
$ws.Range("N93").Value = -5553.6924
$ws.Range("H122").Value = 100000
$ws.Range("I122").Value = 100000
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 300000
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -297550
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 576.2414
$ws.Range("I107").Value = 492.65384
$ws.Range("K107").Value = 1477.96152
$ws.Range("M107").Value = 442.0384799999999
$ws.Range("H109").Value = 13437.5
$ws.Range("J109").Value = 13437.5
$ws.Range("L109").Value = 13437.5
$ws.Range("N109").Value = -16211.5
$ws.Range("H122").Value = 142858850
$ws.Range("I122").Value = 166668320
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 500004960
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -500002510
$ws.Range("N122").Value = -10900
$ws.Range("H132").Value = 840.6667
$ws.Range("I132").Value = 742.7568
$ws.Range("J132").Value = 1021.8
$ws.Range("K132").Value = 2228.2704
$ws.Range("L132").Value = 3065.4
$ws.Range("M132").Value = 301.7296000000001
$ws.Range("N132").Value = -8125.4
$ws.Range("H137").Value = 68712.414
$ws.Range("J137").Value = 68712.414
$ws.Range("L137").Value = 68712.414
$ws.Range("N137").Value = -78912.414
